# Update power/temperature data in "imx8 results" for experiment 6
# (adds a "Mod.2 B predictor" data row to each of the three instance blocks).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Instance 4 block: row 44 ("Mod.2 B predictor") ---
$ws.Range("A44").Value = "Mod.2 B predictor"
$ws.Range("B44").Value = 55.96
$ws.Range("C44").Value = 55.26
$ws.Range("D44").Value = 54.856
$ws.Range("L44").Value = 25.5396
$ws.Range("M44").Value = 24.8362
$ws.Range("N44").Value = 24.57

# --- Instance 5 block: row 55 ("Mod. B predictor") ---
$ws.Range("B55").Value = 56.66
$ws.Range("C55").Value = 56.76
$ws.Range("D55").Value = 56.624
$ws.Range("L55").Value = 23.8944
$ws.Range("M55").Value = 23.6648
$ws.Range("N55").Value = 23.4464

# --- Instance 5 block: row 56 ("No B optimal LTF") ---
$ws.Range("B56").Value = 58.008
$ws.Range("C56").Value = 58.024
$ws.Range("D56").Value = 57.836
$ws.Range("L56").Value = 23.7828
$ws.Range("M56").Value = 23.5702
$ws.Range("N56").Value = 23.356

# --- Instance 5 block: row 58 ("Mod.2 B predictor") ---
$ws.Range("A58").Value = "Mod.2 B predictor"
$ws.Range("B58").Value = 57.146
$ws.Range("C58").Value = 57.264
$ws.Range("D58").Value = 56.964
$ws.Range("L58").Value = 24.452
$ws.Range("M58").Value = 24.2956
$ws.Range("N58").Value = 24.154

# Match the author's final selection before saving.
$ws.Range("N56").Select() | Out-Null
